$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1067
$ws.Range("F8").Value = 577
$ws.Range("F9").Value = 1506
$ws.Range("F11").Value = 1399
$ws.Range("F14").Value = 1704
$ws.Range("F15").Value = 1773
$ws.Range("F16").Value = 825
$ws.Range("F17").Value = 256
$ws.Range("F18").Value = 1432
$ws.Range("F21").Value = 1163
$ws.Range("F22").Value = 379
$ws.Range("F23").Value = 421
$ws.Range("F24").Value = 47
$ws.Range("F25").Value = 3839
$ws.Range("F26").Value = 719
$ws.Range("F28").Value = 1600
$ws.Range("F30").Value = 67

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 16
$ws.Range("F9").Value = 39
$ws.Range("F13").Value = 92

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 802
$ws.Range("F3").Value = 26

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 802
$ws.Range("F4").Value = 26
$ws.Range("F12").Value = 16
$ws.Range("F14").Value = 39
$ws.Range("F16").Value = 1067
$ws.Range("F19").Value = 577
$ws.Range("F20").Value = 1506
$ws.Range("F22").Value = 1399
$ws.Range("F23").Value = 3043
$ws.Range("F25").Value = 1704
$ws.Range("F26").Value = 1773
$ws.Range("F27").Value = 825
$ws.Range("F28").Value = 256
$ws.Range("F29").Value = 1432
$ws.Range("F34").Value = 1163
$ws.Range("F35").Value = 379
$ws.Range("F36").Value = 421
$ws.Range("F37").Value = 47
$ws.Range("F38").Value = 3839
$ws.Range("F39").Value = 719
$ws.Range("F41").Value = 1600
$ws.Range("F42").Value = 92
$ws.Range("F45").Value = 67
